# Data Driving Valid Login Test Script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Data row first, so these shared strings are registered before the header's
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Column B auto-fit width (bestFit) to match the diff's column width
$ws.Columns.Item(2).AutoFit() | Out-Null

# Selection / zoom as captured by the saved view state
$ws.Range("A3").Select()
$excel.ActiveWindow.Zoom = 205
